# Fixed release timeline positions
#
# Slide 2 holds the release timeline: for each milestone there is a
# rounded-rectangle "pill" shape plus a label TextBox sitting directly
# beneath it. Several pills/labels had the wrong horizontal (x) offset;
# this corrects them while leaving their y-offsets and sizes untouched.
#
# NOTE on units: the OOXML stores offsets in EMU (914400 EMU = 1 inch),
# but the PowerPoint COM object model's Shape.Left/.Top are expressed in
# points (1 pt = 12700 EMU) and stored internally as 32-bit floats. A
# naive EMU/12700 division can therefore land one EMU short after the
# float32 round-trip back on save (e.g. 7132320 EMU -> 561.6 pt -> stored
# as a float32 that is a hair under 561.6 -> floors back to 7132319 EMU).
# Emu-ToComPoints nudges the point value up by the smallest amount needed
# so the round trip reproduces the exact target EMU value.

function Emu-ToComPoints {
    param($TargetEmu, $EmuPerPoint)
    $pts = $TargetEmu / $EmuPerPoint
    $step = 0.0000001
    $i = 0
    while ($i -lt 2000) {
        $candidate = $pts + ($i * $step)
        $c32 = [double]([single]$candidate)
        $emu = [math]::Floor($c32 * $EmuPerPoint)
        if ($emu -eq $TargetEmu) {
            return $candidate
        }
        $i = $i + 1
    }
    # Fall back to the naive conversion if no exact nudge is found.
    return $pts
}

$EMU_PER_POINT = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Map of Shapes.Item(index) -> new Left, in EMU (converted to points below).
$moves = @{
    17 = 2103120   # Rounded Rectangle 17 (was 731520)
    18 = 2103120   # TextBox 18 - "PHASE 2a / May" (was 731520)
    19 = 2743200   # Rounded Rectangle 19 (was 2011680)
    20 = 2743200   # TextBox 20 - "JUNE / Jun" (was 2011680)
    21 = 5303520   # Rounded Rectangle 21 (was 6217920)
    22 = 5303520   # TextBox 22 - "SEPT / Sep" (was 6217920)
    23 = 7132320   # Rounded Rectangle 23 (was 8961120)
    24 = 7132320   # TextBox 24 - "NOV / Nov" (was 8961120)
    25 = 10515600  # Rounded Rectangle 25 (was 10789920)
    26 = 10515600  # TextBox 26 - "MAR 2027 / Mar" (was 10789920)
}

foreach ($idx in $moves.Keys) {
    $shape = $s.Shapes.Item($idx)
    $shape.Left = Emu-ToComPoints $moves[$idx] $EMU_PER_POINT
}
